# Fix double-space typo "For  <Month>, the ... slot is unavailable"
# -> "For <Month>, the ... slot is unavailable"
# The affected strings live on the "Declined" worksheet (column H).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Declined")

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count
$startRow = $used.Row
$startCol = $used.Column

for ($r = 0; $r -lt $rows; $r++) {
    for ($c = 0; $c -lt $cols; $c++) {
        $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string] -and $val -like "For  *slot is unavailable*") {
            $cell.Value = $val -replace "^For  ", "For "
        }
    }
}
